# Insert a new weekly record row at row 361, pushing existing rows 361-448
# down to 362-449, then populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("361").Insert()

$ws.Range("A361").Value = 3
$ws.Range("B361").Value = "Femacal de La Calera"
$ws.Range("C361").Value = "Coquimbo"
$ws.Range("D361").Value = 44932
$ws.Range("E361").Value = 5
$ws.Range("F361").Value = 100112012
$ws.Range("G361").Value = "Espinaca"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 170
$ws.Range("K361").Value = 4000
$ws.Range("L361").Value = 4500
$ws.Range("M361").Value = 4235
$ws.Range("N361").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O361").Value = "Provincia de Quillota"
$ws.Range("P361").Value = 1412
$ws.Range("Q361").Value = 3
$ws.Range("R361").Value = "Hortaliza"
